$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a numeric-looking string while keeping it stored
# as text (so "248.72" stays exactly "248.72" and not a float / General
# number re-formatted value), and without leaving a custom number format
# behind on the cell.
function Set-TextValue($sheet, $cellRef, $val) {
    $c = $sheet.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Helper: set a cell to a plain (non numeric-looking) text value.
function Set-PlainValue($sheet, $cellRef, $val) {
    $sheet.Range($cellRef).Value = $val
}

# --- Price (column D) refreshes ---
Set-TextValue $ws 'D2' '248.72'
Set-TextValue $ws 'D3' '21.76'
Set-TextValue $ws 'D4' '5.575'
Set-TextValue $ws 'D5' '0.05661'
Set-TextValue $ws 'D6' '6.440'
Set-TextValue $ws 'D7' '0.8009'
Set-TextValue $ws 'D9' '0.1432'
Set-TextValue $ws 'D10' '0.07269'
Set-TextValue $ws 'D11' '0.03130'
Set-TextValue $ws 'D12' '0.02922'
Set-TextValue $ws 'D13' '0.09278'
Set-TextValue $ws 'D14' '0.001666'
Set-TextValue $ws 'D15' '3.212'
Set-TextValue $ws 'D16' '0.04749'

# --- Rows 17-24: coin ranking reshuffled (symbol list update) ---
Set-PlainValue $ws 'B17' 'TigerCash'
Set-PlainValue $ws 'C17' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws 'D17' '0.006460'
Set-PlainValue $ws 'E17' '16TigerCashTCH'

Set-PlainValue $ws 'B18' 'HotbitToken'
Set-PlainValue $ws 'C18' 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws 'D18' '0.005059'
Set-PlainValue $ws 'E18' '17HotbitTokenHTB'

Set-PlainValue $ws 'B19' 'BitKan'
Set-PlainValue $ws 'C19' 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws 'D19' '0.001056'
Set-PlainValue $ws 'E19' '18BitKanKAN'

Set-PlainValue $ws 'B20' 'NitroEx'
Set-PlainValue $ws 'C20' 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextValue $ws 'D20' '0.0001502'
Set-PlainValue $ws 'E20' '19NitroExNTX'

Set-PlainValue $ws 'B21' 'LEO'
Set-PlainValue $ws 'C21' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws 'D21' '3.973'
Set-PlainValue $ws 'E21' '20LEOLEO'

Set-PlainValue $ws 'B22' 'GateToken'
Set-PlainValue $ws 'C22' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws 'D22' '3.378'
Set-PlainValue $ws 'E22' '21GateTokenGT'

Set-PlainValue $ws 'B23' 'BTSEToken'
Set-PlainValue $ws 'C23' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws 'D23' '2.088'
Set-PlainValue $ws 'E23' '22BTSETokenBTSE'

Set-PlainValue $ws 'B24' 'One'
Set-PlainValue $ws 'C24' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws 'D24' '0.01164'
Set-PlainValue $ws 'E24' '23OneONEBestin24h'

# --- More price refreshes ---
Set-TextValue $ws 'D25' '0.3266'
Set-TextValue $ws 'D27' '0.0003204'
Set-TextValue $ws 'D40' '0.04122'
Set-TextValue $ws 'D41' '0.006918'

# --- Rows 42-43: CEJI / BKEXToken swap positions ---
Set-PlainValue $ws 'B42' 'CEJI'
Set-PlainValue $ws 'C42' 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws 'D42' '0.003504'
Set-PlainValue $ws 'E42' '41CEJICEJI'

Set-PlainValue $ws 'B43' 'BKEXToken'
Set-PlainValue $ws 'C43' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws 'D43' '0.1041'
Set-PlainValue $ws 'E43' '42BKEXTokenBKK'

# --- More price refreshes ---
Set-TextValue $ws 'D44' '0.008517'
Set-TextValue $ws 'D45' '0.00005638'

# --- Rows 47-48: Best/Worst in 24h labels move ---
Set-TextValue $ws 'D47' '0.7861'
Set-PlainValue $ws 'E47' '46CoinbaseStockTokenCOIN'

Set-TextValue $ws 'D48' '0.01674'
Set-PlainValue $ws 'E48' '47BOLOBOLOWorstin24h'

Set-TextValue $ws 'D49' '0.00002102'
